# Fix Training Data Issue (#48)
# The per-game team box-score data in this sheet was pulled one day off
# (2007-08 season, NBA game date 2008-02-03 mislabeled as "2-3-2007-08").
# This patches the Date column (BF) to the correct ISO date string and
# corrects the handful of ranked/aggregate stat cells that shifted as a
# result of pulling the correct day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column (BF2:BF31): "2-3-2007-08" -> "2008-02-03" ---
# Force text format first so Excel does not reinterpret the ISO-looking
# string as a date serial value.
$ws.Range("BF2:BF31").NumberFormat = "@"
$ws.Range("BF2").Value = "2008-02-03"
$ws.Range("BF3").Value = "2008-02-03"
$ws.Range("BF4").Value = "2008-02-03"
$ws.Range("BF5").Value = "2008-02-03"
$ws.Range("BF6").Value = "2008-02-03"
$ws.Range("BF7").Value = "2008-02-03"
$ws.Range("BF8").Value = "2008-02-03"
$ws.Range("BF9").Value = "2008-02-03"
$ws.Range("BF10").Value = "2008-02-03"
$ws.Range("BF11").Value = "2008-02-03"
$ws.Range("BF12").Value = "2008-02-03"
$ws.Range("BF13").Value = "2008-02-03"
$ws.Range("BF14").Value = "2008-02-03"
$ws.Range("BF15").Value = "2008-02-03"
$ws.Range("BF16").Value = "2008-02-03"
$ws.Range("BF17").Value = "2008-02-03"
$ws.Range("BF18").Value = "2008-02-03"
$ws.Range("BF19").Value = "2008-02-03"
$ws.Range("BF20").Value = "2008-02-03"
$ws.Range("BF21").Value = "2008-02-03"
$ws.Range("BF22").Value = "2008-02-03"
$ws.Range("BF23").Value = "2008-02-03"
$ws.Range("BF24").Value = "2008-02-03"
$ws.Range("BF25").Value = "2008-02-03"
$ws.Range("BF26").Value = "2008-02-03"
$ws.Range("BF27").Value = "2008-02-03"
$ws.Range("BF28").Value = "2008-02-03"
$ws.Range("BF29").Value = "2008-02-03"
$ws.Range("BF30").Value = "2008-02-03"
$ws.Range("BF31").Value = "2008-02-03"

# --- Corrected stat values, row by row ---
# Row 3
$ws.Range("AK3").Value = 5
# Row 4
$ws.Range("AJ4").Value = 19
$ws.Range("AS4").Value = 25
# Row 5
$ws.Range("AD5").Value = 14
# Row 6
$ws.Range("AD6").Value = 14
$ws.Range("AT6").Value = 4
$ws.Range("BA6").Value = 22
# Row 7
$ws.Range("D7").Value = 45
$ws.Range("F7").Value = 14
$ws.Range("G7").Value = 0.6889999999999999
$ws.Range("I7").Value = 36.8
$ws.Range("K7").Value = 0.471
$ws.Range("L7").Value = 5.9
$ws.Range("N7").Value = 0.353
$ws.Range("Q7").Value = 0.823
$ws.Range("U7").Value = 20.3
$ws.Range("V7").Value = 13
$ws.Range("Z7").Value = 22.2
$ws.Range("AA7").Value = 21.8
$ws.Range("AB7").Value = 100.9
$ws.Range("AC7").Value = 4.9
$ws.Range("AD7").Value = 24
$ws.Range("AF7").Value = 3
$ws.Range("AH7").Value = 24
$ws.Range("AI7").Value = 11
$ws.Range("AK7").Value = 4
$ws.Range("AL7").Value = 18
$ws.Range("AO7").Value = 6
$ws.Range("AP7").Value = 10
$ws.Range("AU7").Value = 21
$ws.Range("AX7").Value = 5
$ws.Range("AZ7").Value = 23
$ws.Range("BB7").Value = 8
$ws.Range("BC7").Value = 7
# Row 8
$ws.Range("AD8").Value = 14
# Row 9
$ws.Range("D9").Value = 46
$ws.Range("E9").Value = 33
$ws.Range("G9").Value = 0.717
$ws.Range("L9").Value = 5.8
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 0.366
$ws.Range("O9").Value = 18.8
$ws.Range("P9").Value = 24.7
$ws.Range("Q9").Value = 0.761
$ws.Range("S9").Value = 29.3
$ws.Range("T9").Value = 40.8
$ws.Range("V9").Value = 11.7
$ws.Range("X9").Value = 5.1
$ws.Range("Z9").Value = 20.5
$ws.Range("AA9").Value = 20.3
$ws.Range("AB9").Value = 97.7
$ws.Range("AC9").Value = 7.3
$ws.Range("AD9").Value = 14
$ws.Range("AI9").Value = 14
$ws.Range("AJ9").Value = 20
$ws.Range("AL9").Value = 19
$ws.Range("AN9").Value = 12
$ws.Range("AO9").Value = 13
$ws.Range("AP9").Value = 17
$ws.Range("AS9").Value = 26
$ws.Range("AW9").Value = 15
$ws.Range("AX9").Value = 11
$ws.Range("BA9").Value = 21
# Row 10
$ws.Range("AE10").Value = 8
# Row 11
$ws.Range("AI11").Value = 16
$ws.Range("AW11").Value = 12
# Row 12
$ws.Range("AN12").Value = 10
# Row 14
$ws.Range("D14").Value = 45
$ws.Range("E14").Value = 29
$ws.Range("G14").Value = 0.644
$ws.Range("I14").Value = 39
$ws.Range("J14").Value = 82
$ws.Range("K14").Value = 0.476
$ws.Range("N14").Value = 0.369
$ws.Range("O14").Value = 21.9
$ws.Range("P14").Value = 29.2
$ws.Range("Q14").Value = 0.749
$ws.Range("R14").Value = 11
$ws.Range("S14").Value = 33.8
$ws.Range("T14").Value = 44.7
$ws.Range("U14").Value = 23.5
$ws.Range("W14").Value = 8.199999999999999
$ws.Range("X14").Value = 5.2
$ws.Range("AA14").Value = 23.1
$ws.Range("AB14").Value = 107.2
$ws.Range("AC14").Value = 5.7
$ws.Range("AD14").Value = 24
$ws.Range("AE14").Value = 8
$ws.Range("AH14").Value = 24
$ws.Range("AT14").Value = 3
$ws.Range("AV14").Value = 24
$ws.Range("AW14").Value = 5
$ws.Range("AX14").Value = 9
$ws.Range("BA14").Value = 5
$ws.Range("BC14").Value = 6
# Row 15
$ws.Range("AP15").Value = 19
$ws.Range("AU15").Value = 22
$ws.Range("AX15").Value = 6
$ws.Range("BB15").Value = 9
# Row 16
$ws.Range("AD16").Value = 14
# Row 17
$ws.Range("BA17").Value = 23
# Row 18
$ws.Range("AD18").Value = 14
$ws.Range("AI18").Value = 13
# Row 19
$ws.Range("AS19").Value = 17
$ws.Range("BA19").Value = 6
# Row 20
$ws.Range("AD20").Value = 14
$ws.Range("AW20").Value = 16
# Row 21
$ws.Range("AO21").Value = 14
$ws.Range("AP21").Value = 11
# Row 22
$ws.Range("AI22").Value = 12
$ws.Range("AK22").Value = 6
$ws.Range("AN22").Value = 11
# Row 23
$ws.Range("AW23").Value = 6
# Row 24
$ws.Range("AE24").Value = 2
$ws.Range("AH24").Value = 28
# Row 25
$ws.Range("AD25").Value = 14
$ws.Range("BC25").Value = 15
# Row 26
$ws.Range("AD26").Value = 14
$ws.Range("AQ26").Value = 4
$ws.Range("BB26").Value = 10
# Row 27
$ws.Range("AD27").Value = 24
$ws.Range("AE27").Value = 8
$ws.Range("AG27").Value = 6
$ws.Range("BC27").Value = 8
# Row 29
$ws.Range("AD29").Value = 14
$ws.Range("AJ29").Value = 12
$ws.Range("BB29").Value = 14
# Row 30
$ws.Range("AO30").Value = 7
$ws.Range("BC30").Value = 5
# Row 31
$ws.Range("D31").Value = 45
$ws.Range("F31").Value = 21
$ws.Range("G31").Value = 0.533
$ws.Range("H31").Value = 48.6
$ws.Range("I31").Value = 36.3
$ws.Range("L31").Value = 6.6
$ws.Range("O31").Value = 19.5
$ws.Range("P31").Value = 24.6
$ws.Range("Q31").Value = 0.795
$ws.Range("S31").Value = 30.4
$ws.Range("T31").Value = 42.6
$ws.Range("V31").Value = 14
$ws.Range("W31").Value = 7.7
$ws.Range("AB31").Value = 98.7
$ws.Range("AC31").Value = 0.9
$ws.Range("AD31").Value = 24
$ws.Range("AF31").Value = 15
$ws.Range("AI31").Value = 15
$ws.Range("AJ31").Value = 11
$ws.Range("AP31").Value = 18
$ws.Range("AQ31").Value = 3
$ws.Range("AS31").Value = 16
$ws.Range("AW31").Value = 11
$ws.Range("AX31").Value = 10
$ws.Range("BB31").Value = 13
$ws.Range("BC31").Value = 14
